$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("universal")

# Row 3: new "teste" / "pediatrica" (pediatric dosing) entry
$ws.Range("A3").Value = "teste"
$ws.Range("B3").Value = "pediatrica"
$ws.Range("C3").Value = "sim"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 5110
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = "mg/kg"
$ws.Range("L3").Value = 3
$ws.Range("M3").Value = "12, 8"
$ws.Range("N3").Value = "peso, idade, dosagem, dose, intervalo, via"

# Row 4: new "teste" / "adulta" (adult dosing) entry
$ws.Range("A4").Value = "teste"
$ws.Range("B4").Value = "adulta"
$ws.Range("D4").Value = 5475
$ws.Range("E4").Value = 54750
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 6
$ws.Range("K4").Value = "kg"
$ws.Range("L4").Value = 6
$ws.Range("N4").Value = "idade, peso, dosagem, dose"

# Column N ("campos") was widened to fit the new, longer text entries
$ws.Columns.Item(14).ColumnWidth = 36.2916667

# Final selection left on the last-edited cell
$ws.Range("N4").Select()
